$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove all existing comments on row 15 (they do not auto-shift on column insert) ---
$origCommentCols = @('A', 'B', 'C', 'D', 'E', 'F', 'G', 'H', 'I', 'J', 'K', 'L', 'M', 'N', 'O', 'P', 'Q', 'R', 'S', 'T', 'U', 'V', 'W', 'X', 'Y', 'Z', 'AA', 'AB', 'AC', 'AD', 'AE', 'AF', 'AG', 'AH', 'AI', 'AJ', 'AM', 'AN', 'AO', 'AQ', 'AR', 'AS', 'AT', 'AU', 'AV', 'AW', 'AX', 'AY', 'AZ', 'BA', 'BB', 'BC', 'BD', 'BE', 'BF', 'BG', 'BH', 'BI', 'BJ', 'BK', 'BL', 'BM', 'BN', 'BO', 'BP', 'BQ', 'BR', 'BS', 'BT', 'BU', 'BV', 'BW', 'BX', 'BY')
foreach ($col in $origCommentCols) {
    $ref = $col + "15"
    $cm = $ws.Range($ref).Comment
    if ($cm -ne $null) {
        $cm.Delete()
    }
}

# --- Step 2: insert the two new columns at their correct positions ---
# "culture_collection" goes where "douche" used to be (column Z), pushing douche and
# everything after it one column to the right.
$ws.Columns("Z").Insert()
# "specimen_voucher" goes where "temperature" used to be. After the first insertion,
# "temperature" now sits at column BX, so insert there.
$ws.Columns("BX").Insert()

# --- Step 3: set header text for the two new columns ---
$ws.Range("Z15").Value = "culture_collection"
$ws.Range("BX15").Value = "specimen_voucher"

# --- Step 4: (re)apply every comment on row 15 at its correct final position ---
$finalComments = @(
    @{ col = 'A'; text = 'The sample name is a name that you choose for the sample, it works as an ID.  Each sample name must be unique in samples submitted from a submission account.' },
    @{ col = 'B'; text = 'Sample title should be short and informative. Each sample title must be unique in a submission.  Examples: 1) Escherichia coli O104:H4 str. C227-11 clinical isolate 2010_333_NC-6;  2) CD8+ T cells from female TSG6-knockout BALB/c mouse;  3) Human metagenome isolated from urine of healthy female.' },
    @{ col = 'C'; text = 'A brief description for the sample.' },
    @{ col = 'D'; text = 'The most descriptive organism name for this sample (to the species, if relevant) in the NCBI Taxonomy database, http://www.ncbi.nlm.nih.gov/taxonomy If it is not in the database, provide as much information about the organism as possible and the DDBJ staff apply a new organism name to NCBI Taxonomy.' },
    @{ col = 'E'; text = 'NCBI Taxonomy identifier. This is appropriate for individual organisms, some metagenomes and environmental samples (http://www.ncbi.nlm.nih.gov/Taxonomy/Browser/wwwtax.cgi?mode=Undef&id=12908&lvl=3&lin=f&keep=1&srchmode=1&unlock).  If it is not in the database, enter a tentative ID (e.g., 1). The DDBJ staff apply a new organism name to NCBI Taxonomy, and then the tentative ID is replaced by an assigned TaxID.' },
    @{ col = 'F'; text = 'Associated BioProject accession number (PRJDB)' },
    @{ col = 'G'; text = 'Organism group

microbial or eukaryotic strain name' },
    @{ col = 'H'; text = 'Organism group

Identification or description of the specific individual from which this sample was obtained' },
    @{ col = 'I'; text = 'Organism group

Cultivar name - cultivated variety of plant' },
    @{ col = 'J'; text = 'Organism group

a population within a given species displaying genetically based, phenotypic traits that reflect adaptation to a local habitat, e.g., Columbia' },
    @{ col = 'K'; text = 'Time of sampling (single instance or interval, eg., 2008-01-23T19:23:10, 2008-01-23, 2008-01, 2008, 1952-10-21T11:43Z/1952-10-21T17:43Z, 1952-10-21/1953-02-15, 1952-10/1953-02, 1952/1953)' },
    @{ col = 'L'; text = 'Descriptor of the broad ecological context of a sample. Examples include: desert, taiga or deciduous woodland. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO' },
    @{ col = 'M'; text = 'Descriptor of the local environment. Examples include: harbor, cliff, or lake. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO' },
    @{ col = 'N'; text = 'Material that was displaced by the sample, or material in which a sample was embedded, prior to the sampling event. Examples include: air, soil, or water. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO' },
    @{ col = 'O'; text = 'Estimated size of genome' },
    @{ col = 'P'; text = 'Geographical origin of the sample; use the appropriate name from the list, http://www.ddbj.nig.ac.jp/sub/country-e.html. Use a colon to separate the country or ocean from more detailed information about the location, eg "Japan:Kanagawa, Hakone, Lake Ashi" ' },
    @{ col = 'Q'; text = 'The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".' },
    @{ col = 'R'; text = 'Publication reference in the form of pubmed ID, DOI or URL for isolation and growth condition specifications of the organism/material' },
    @{ col = 'S'; text = 'The geographical coordinates of the location where the sample was collected. Specify as decimal degrees latitude and longitude in format "d[d.dddd] N|S d[dd.dddd] W|E", eg, 47.94 N 28.12 W' },
    @{ col = 'T'; text = 'Reports the number of replicons in a nuclear genome of eukaryotes, in the genome of a bacterium or archaea or the number of segments in a segmented virus. Always applied to the haploid chromosome count of a eukaryote' },
    @{ col = 'U'; text = 'The ploidy level of the genome (e.g. allopolyploid, haploid, diploid, triploid, tetraploid). ' },
    @{ col = 'V'; text = 'This field is specific to different taxa. For phage: lytic/lysogenic/temperate/obligately lytic;  for plasmid: incompatibility group;  for eukaryote: asexual/sexual' },
    @{ col = 'W'; text = 'Free-living or from host (define relationship)' },
    @{ col = 'X'; text = 'specification of birth control medication used' },
    @{ col = 'Y'; text = 'list of chemical compounds administered to the host or site where sampling occurred, and when (e.g. antibiotics, N fertilizer, air filter); can include multiple compounds. For Chemical Entities of Biological Interest ontology (CHEBI) (v1.72), please see http://bioportal.bioontology.org/visualize/44603' },
    @{ col = 'Z'; text = 'Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier' },
    @{ col = 'AA'; text = 'date of most recent douche' },
    @{ col = 'AB'; text = 'ethnicity of the subject' },
    @{ col = 'AC'; text = 'Plasmids that have significance phenotypic consequence' },
    @{ col = 'AD'; text = 'history of gynecological disorders; can include multiple disorders' },
    @{ col = 'AE'; text = 'Health or disease status of sample at time of collection' },
    @{ col = 'AF'; text = 'Age of host at the time of sampling' },
    @{ col = 'AG'; text = 'body mass index of the host, calculated as weight/(height)squared' },
    @{ col = 'AH'; text = 'substance produced by the host, e.g. stool, mucus, where the sample was obtained from' },
    @{ col = 'AI'; text = 'core body temperature of the host when sample was collected' },
    @{ col = 'AJ'; text = 'type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types' },
    @{ col = 'AK'; text = 'Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh' },
    @{ col = 'AN'; text = 'the height of subject' },
    @{ col = 'AO'; text = 'content of last meal and time since feeding; can include multiple values' },
    @{ col = 'AP'; text = 'most frequent job performed by subject' },
    @{ col = 'AR'; text = 'resting pulse of the host, measured as beats per minute' },
    @{ col = 'AS'; text = 'Gender or physical sex of the host' },
    @{ col = 'AT'; text = 'a unique identifier by which each subject can be referred to, de-identified, e.g. #131' },
    @{ col = 'AU'; text = 'NCBI taxonomy ID of the host, e.g. 9606' },
    @{ col = 'AV'; text = 'Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005' },
    @{ col = 'AW'; text = 'total mass of the host at collection, the unit depends on host' },
    @{ col = 'AX'; text = 'whether subject had hormone replacement theraphy, and if yes start date' },
    @{ col = 'AY'; text = 'specification of whether hysterectomy was performed' },
    @{ col = 'AZ'; text = 'can include multiple medication codes' },
    @{ col = 'BA'; text = 'Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.' },
    @{ col = 'BB'; text = 'A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html' },
    @{ col = 'BC'; text = 'whether full medical history was collected' },
    @{ col = 'BD'; text = 'date of most recent menstruation' },
    @{ col = 'BE'; text = 'date of onset of menopause' },
    @{ col = 'BF'; text = 'any other measurement performed or parameter collected, that is not listed here' },
    @{ col = 'BG'; text = 'total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts' },
    @{ col = 'BH'; text = 'oxygenation status of sample' },
    @{ col = 'BI'; text = 'To what is the entity pathogenic' },
    @{ col = 'BJ'; text = 'type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types' },
    @{ col = 'BK'; text = 'date due of pregnancy' },
    @{ col = 'BL'; text = 'Primary publication or genome report in the form of pubmed ID, DOI or URL' },
    @{ col = 'BM'; text = 'Method or device employed for collecting sample' },
    @{ col = 'BN'; text = 'Processing applied to the sample during or after isolation' },
    @{ col = 'BO'; text = 'salinity of sample, i.e. measure of total salt concentration' },
    @{ col = 'BP'; text = 'Amount or size of sample (volume, mass or area) that was collected' },
    @{ col = 'BQ'; text = 'duration for which sample was stored' },
    @{ col = 'BR'; text = 'location at which sample was stored, usually name of a specific freezer/room' },
    @{ col = 'BS'; text = 'temperature at which sample was stored, e.g. -80' },
    @{ col = 'BT'; text = 'volume (mL) or weight (g) of sample processed for DNA extraction' },
    @{ col = 'BU'; text = 'current sexual partner and frequency of sex' },
    @{ col = 'BV'; text = 'unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.' },
    @{ col = 'BW'; text = 'Information about the genetic distinctness of the lineage (eg., biovar, serovar)' },
    @{ col = 'BX'; text = 'Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier' },
    @{ col = 'BY'; text = 'temperature of the sample at time of sampling' },
    @{ col = 'BZ'; text = 'Feeding position in food chain (eg., chemolithotroph)' },
    @{ col = 'CA'; text = 'history of urogenital disorders, can include multiple disorders' }
)
foreach ($entry in $finalComments) {
    $ref = $entry.col + "15"
    $ws.Range($ref).AddComment($entry.text)
}

Write-Output "done"